# Architecture diagram clean-up: the "Events Center" node (and the
# elbow connector hanging off of it) is removed from the diagram; the
# remaining "Logs Center" node/connector are left exactly as they were.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Remove the dangling connector that originates at the "Oval 5" shape
# (id=6 / "Events Center") before removing the shape itself.
$s.Shapes.Item("Elbow Connector 8").Delete()

# Remove the "Events Center" oval.
$s.Shapes.Item("Oval 5").Delete()
